$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Hoja1" -> "Variable_problema2" and populate it with the new
#    exercises (Bloque 4 - descomposicion de numeros en unidades/decenas/...).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Hoja1")
$ws4.Name = "Variable_problema2"

# --- Row 1 ---
$ws4.Range("A1").Value = 1
$ws4.Range("B1").Formula = "=MOD(A1,10)"
$ws4.Range("C1").Formula = "=ROUNDDOWN(A1/10,0)"

# --- Row 2 ---
$ws4.Range("A2").Value = 12
$ws4.Range("B2").Formula = "=MOD(A2,10)"
$ws4.Range("C2").Formula = "=ROUNDDOWN(A2/10,0)"
$ws4.Range("D2").Formula = "=MOD(C2,10)"
$ws4.Range("E2").Formula = "=ROUNDDOWN(C2/10,0)"

# --- Row 3 ---
$ws4.Range("A3").Value = 105
$ws4.Range("B3").Formula = "=MOD(A3,10)"
$ws4.Range("C3").Formula = "=ROUNDDOWN(A3/10,0)"
$ws4.Range("D3").Formula = "=MOD(C3,10)"
$ws4.Range("E3").Formula = "=ROUNDDOWN(C3/10,0)"
$ws4.Range("F3").Formula = "=MOD(E3,10)"
$ws4.Range("G3").Formula = "=ROUNDDOWN(E3/10,0)"

# --- Row 7 ---
$ws4.Range("D7").Value = 4
$ws4.Range("E7").Value = 3
$ws4.Range("F7").Value = 2
$ws4.Range("G7").Value = 1

# --- Row 8 ---
$ws4.Range("A8").Value = 4567
$ws4.Range("B8").Value = 4
$ws4.Range("D8").Formula = "=ROUNDDOWN(MOD(`$A`$8,POWER(10,4)),0)"
$ws4.Range("E8").Formula = "=ROUNDDOWN(MOD(`$A`$8,POWER(10,2-1)),0)"
$ws4.Range("F8").Formula = "=ROUNDDOWN(MOD(`$A`$8,POWER(10,4-1)),0)"
$ws4.Range("G8").Formula = "=ROUNDDOWN(MOD(`$A`$8,POWER(10,4-1)),0)"

# --- Row 11 ---
$ws4.Range("D11").Formula = "=ROUNDDOWN(`$A`$8/POWER(10,4-1),0)"
$ws4.Range("E11").Formula = "=ROUNDDOWN(MOD(D11,10),0)"
$ws4.Range("F11").Formula = "=ROUNDDOWN(`$A`$8/POWER(10,3-1),0)"
$ws4.Range("G11").Formula = "=ROUNDDOWN(MOD(F11,10),0)"
$ws4.Range("H11").Formula = "=ROUNDDOWN(`$A`$8/POWER(10,2-1),0)"
$ws4.Range("I11").Formula = "=ROUNDDOWN(MOD(H11,10),0)"
$ws4.Range("J11").Formula = "=`$A`$8"
$ws4.Range("K11").Formula = "=ROUNDDOWN(MOD(J11,10),0)"

# --- Row 13 ---
$ws4.Range("E13").Formula = "=4*1"
$ws4.Range("G13").Formula = "=4+5*10"
$ws4.Range("I13").Formula = "=54+6*100"
$ws4.Range("K13").Formula = "=654+7*1000"

# Make this the active sheet / selected cell (matches tabSelected + selection
# in the saved file) and restore the tab order's active index.
$ws4.Activate()
$ws4.Range("E13").Select()

# ---------------------------------------------------------------------------
# 2. PS3_14: tidy up the two shared-formula groups whose declared "ref" span
#    had grown one cell too wide (E2:K2 -> F2:K2, F19:S19 -> F19:R19). Excel
#    recomputes these spans to exactly match the cells actually sharing the
#    formula; re-asserting the same formulas forces that recompute.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PS3_14")

$ws3.Range("F2").Formula = "=E2+POWER(-1,F1+1)*POWER(`$A`$2,2*F1-1)/FACT(2*F1-1)"
$ws3.Range("G2").Formula = "=F2+POWER(-1,G1+1)*POWER(`$A`$2,2*G1-1)/FACT(2*G1-1)"
$ws3.Range("H2").Formula = "=G2+POWER(-1,H1+1)*POWER(`$A`$2,2*H1-1)/FACT(2*H1-1)"
$ws3.Range("I2").Formula = "=H2+POWER(-1,I1+1)*POWER(`$A`$2,2*I1-1)/FACT(2*I1-1)"
$ws3.Range("J2").Formula = "=I2+POWER(-1,J1+1)*POWER(`$A`$2,2*J1-1)/FACT(2*J1-1)"
$ws3.Range("K2").Formula = "=J2+POWER(-1,K1+1)*POWER(`$A`$2,2*K1-1)/FACT(2*K1-1)"

$ws3.Range("F19").Formula = "=ABS(F12)-ABS(E12)"
$ws3.Range("G19").Formula = "=ABS(G12)-ABS(F12)"
$ws3.Range("H19").Formula = "=ABS(H12)-ABS(G12)"
$ws3.Range("I19").Formula = "=ABS(I12)-ABS(H12)"
$ws3.Range("J19").Formula = "=ABS(J12)-ABS(I12)"
$ws3.Range("K19").Formula = "=ABS(K12)-ABS(J12)"
$ws3.Range("L19").Formula = "=ABS(L12)-ABS(K12)"
$ws3.Range("M19").Formula = "=ABS(M12)-ABS(L12)"
$ws3.Range("N19").Formula = "=ABS(N12)-ABS(M12)"
$ws3.Range("O19").Formula = "=ABS(O12)-ABS(N12)"
$ws3.Range("P19").Formula = "=ABS(P12)-ABS(O12)"
$ws3.Range("Q19").Formula = "=ABS(Q12)-ABS(P12)"
$ws3.Range("R19").Formula = "=ABS(R12)-ABS(Q12)"

# Restore PS3_14's own selection / scroll anchor.
$ws3.Range("E13").Select()
try { $excel.ActiveWindow.ScrollRow = 4 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

# Leave the workbook with Variable_problema2 as the active/selected tab.
$ws4.Activate()
$ws4.Range("E13").Select()
